# BGU-105 WPF XAML generator (from structs): add a new sheet "WPF_temp_test"
# that generates the boilerplate C# line
#   quCtrl.Content = new <ClassName>(); MessageBox.Show("Press OK to continue...");
# for every "Info"-like class name used elsewhere in the workbook.

$wb = $excel.ActiveWorkbook

# Alphabetically (case-insensitive) sorted list of class names gathered from
# the rest of the workbook (ClassName-ish strings already present as shared
# strings elsewhere in the file).
$classNames = @(
    "BankInfo",
    "BankruptcyInvestigationInfo",
    "BreachOfLawRecordInfo",
    "ContactInfo",
    "CountryInfo",
    "CourtDecisionInfo",
    "CourtInfo",
    "CreditRatingInfo",
    "CurrencyAmount",
    "EconomicActivityType",
    "EmailInfo",
    "FinancialOversightAuthorityInfo",
    "GenericPersonID",
    "GenericPersonInfo",
    "ImperfectBusinessReputationInfo",
    "IncomeOriginInfo",
    "IndebtnessInfo",
    "IPOSharesPurchaseInfo",
    "LegalPersonInfo",
    "LegalTransactionInfo",
    "LiquidatedEntityOwnershipInfo",
    "LoanInfo",
    "LocationInfo",
    "LPRegisteredDateRecordId",
    "MissingInformationResonInfo",
    "OwnershipStructure",
    "OwnershipSummaryInfo",
    "OwnershipVotesInfo",
    "PaymentDeadlineInfo",
    "PersonsAssociation",
    "PhoneInfo",
    "PhysicalPersonInfo",
    "PowerOfAttorneyInfo",
    "PowerOfAttorneySharesPurchaseInfo",
    "RatingAgencyInfo",
    "RegistrarAuthority",
    "RegLicAppx2OwnershipAcqRequestLP",
    "SecondaryMarketSharesPurchaseInfo",
    "SignatoryInfo",
    "SignificantOrDecisiveInfulenceInfo",
    "SignificantOwnershipAcquisitionWaysInfo",
    "TotalOwnershipDetailsInfo",
    "TotalOwnershipSummaryInfo"
)

# Add the new worksheet as the last (7th) tab.
$ws = $wb.Worksheets.Add()
$ws.Name = "WPF_temp_test"
$target = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Move($null, $target)
# NB: `Move` rebinds by position in this host, so re-fetch the sheet object
# by name afterwards instead of continuing to use the stale `$ws` reference.
$ws = $wb.Worksheets("WPF_temp_test")

# Column A (class names) was typed in first, so the one brand-new shared
# string it introduces ("MissingInformationResonInfo") lands in the shared
# string table ahead of the header strings typed afterwards.
$row = 2
foreach ($name in $classNames) {
    $ws.Cells.Item($row, 1).Value = $name
    $row = $row + 1
}

$lastRow = $row - 1

# Header row, typed after the data column.
$ws.Range("A1").Value = "TypeName"
$ws.Range("B1").Value = "quCtrl.Content = new "
$ws.Range("C1").Value = "(); MessageBox.Show(`"Press OK to continue...`");"

# Formula column, filled last (fill-down from B2).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("B$r").Formula = "=`$B`$1&A$r&`$C`$1"
}

# Selection / view state matching the author's final snapshot.
$ws.Range("B2:B$lastRow").Select()
$excel.ActiveWindow.ScrollRow = 19

# Make this sheet the active one (tabSelected) and bump the workbook's
# remembered active tab index.
$ws.Activate()
